$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.210.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5136"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2598"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06454"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.302"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.650.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8038"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.212.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.418"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.01%  "
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.026"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.804"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.77%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1178"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.016"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05113"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.242"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.366"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.232"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.563"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.734"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9255"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.353"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5737"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.166.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01590"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "mCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.555"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.715"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8246"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.799.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4544"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.876"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.60%  "
